$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the header in C1 from 'gender' to 'sex'
$ws.Range("C1").Value = "sex"
